# Weekly fruit/vegetable price update.
# Two new observation rows are inserted into the price history table:
#   - a new row at sheet row 11 (most recent date, 44838)
#   - a new row at sheet row 16 (date 44168)
# Inserting shifts all subsequent rows down, which reproduces the
# observed "row N now holds what used to be row N-1" pattern in the
# rest of the table, and grows the used range from A1:R36 to A1:R38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the first new row at position 11 (pushes old rows 11-36 down to 12-37)
$ws.Rows.Item(11).Insert()

# Insert the second new row at position 16 (pushes rows 16-37 down to 17-38)
$ws.Rows.Item(16).Insert()

# Fill in the brand-new row 11
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44838
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = 100112045
$ws.Cells.Item(11, 7).Value = "Zapallo"
$ws.Cells.Item(11, 8).Value = "Camote"
$ws.Cells.Item(11, 9).Value = "1a nueva(o)"
$ws.Cells.Item(11, 10).Value = 800
$ws.Cells.Item(11, 11).Value = 900
$ws.Cells.Item(11, 12).Value = 950
$ws.Cells.Item(11, 13).Value = 925
$ws.Cells.Item(11, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(11, 15).Value = "Perú"
$ws.Cells.Item(11, 16).Value = 925
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# Fill in the brand-new row 16
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 44168
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100112045
$ws.Cells.Item(16, 7).Value = "Zapallo"
$ws.Cells.Item(16, 8).Value = "Camote"
$ws.Cells.Item(16, 9).Value = "1a nueva(o)"
$ws.Cells.Item(16, 10).Value = 1200
$ws.Cells.Item(16, 11).Value = 1500
$ws.Cells.Item(16, 12).Value = 1700
$ws.Cells.Item(16, 13).Value = 1600
$ws.Cells.Item(16, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 1600
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = "Hortaliza"
